$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 24, shifting existing rows 24-47 down to 25-48
$ws.Rows.Item(24).Insert()

# Populate the new row 24 with data
$ws.Cells.Item(24, 1).Value = 2
$ws.Cells.Item(24, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(24, 3).Value = "Coquimbo"
$ws.Cells.Item(24, 4).Value = 44574
$ws.Cells.Item(24, 5).Value = 4
$ws.Cells.Item(24, 6).Value = 100112032
$ws.Cells.Item(24, 7).Value = "Zapallo italiano"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 340
$ws.Cells.Item(24, 11).Value = 11000
$ws.Cells.Item(24, 12).Value = 12000
$ws.Cells.Item(24, 13).Value = 11500
$ws.Cells.Item(24, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(24, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(24, 16).Value = 192
$ws.Cells.Item(24, 17).Value = 60
$ws.Cells.Item(24, 18).Value = "Hortaliza"
